# Add a new "2022-Q1" quarterly detail sheet (positioned right before the
# existing "总计" sheet) and prepend a corresponding summary row on the
# "总计" sheet.
#
# NOTE: worksheet object references in this host resolve by position, not
# sticky identity -- once a sheet is inserted/moved, any variable that
# pointed at a sheet occupying the shifted slot now refers to whatever
# sheet lives there. So sheets are (re-)looked up by name right before
# each use below, instead of caching references across structural edits.
#
# NOTE: Range.Copy() + Range.PasteSpecial() with no arguments does not
# actually transfer cell formatting (font/border) in this host -- only
# PasteSpecial(xlPasteFormats) / PasteSpecial(-4122) reliably does, so
# formatting is copied in its own explicit pass, separate from values.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Pull formatting (bold header/index font + thin borders + centering) from
# the "2021-Q4" sheet, which has the same column layout (A..H) that
# 2022-Q1 needs. Copy in 3 passes so every one of the 35 rows needed gets
# the correct look (2021-Q4 itself only has 21 rows to copy from).
$templateSheet.Range("A1:H21").Copy()
$newSheet.Range("A1").PasteSpecial($xlPasteFormats)

$templateSheet.Range("A2:H14").Copy()
$newSheet.Range("A22:H34").PasteSpecial($xlPasteFormats)

$templateSheet.Range("A21:H21").Copy()
$newSheet.Range("A35:H35").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# 3. Fund holding rows (2..34). Columns D/E/F/G are stored as text in the
#    source data (e.g. "28.78"), so force text format before writing them.
# ---------------------------------------------------------------------
$rows = @(
    @("008131", "景顺长城竞争优势混合", "28.78", "89.55", "4.54", "1.3066", 7),
    @("005847", "富国沪港深业绩驱动混合A", "44.03", "74.38", "2.95", "1.2989", 6),
    @("001605", "富兰克林国海沪港深成长精选股票", "39.43", "86.45", "2.81", "1.1080", 6),
    @("010104", "景顺长城消费精选混合A", "18.38", "90.93", "4.42", "0.8124", 7),
    @("000934", "国富大中华精选混合QDII", "25.71", "83.59", "3.07", "0.7893", 6),
    @("006370", "国富大中华精选混合QDII美元", "25.71", "83.59", "3.07", "0.7893", 6),
    @("008712", "景顺长城品质成长混合", "15.06", "90.75", "4.42", "0.6657", 7),
    @("160526", "博时优势企业3年封闭运作灵活配置混合（LOF）A", "14.81", "82.86", "4.38", "0.6487", 6),
    @("009983", "永赢港股通品质生活慧选混合", "9.75", "87.97", "6.53", "0.6367", 3),
    @("009846", "富兰克林国海港股通远见价值混合", "19.47", "86.72", "3.07", "0.5977", 4),
    @("005644", "广发沪港深行业龙头混合", "13.85", "88.27", "3.83", "0.5305", 9),
    @("001583", "安信新常态沪港深精选股票", "8.85", "89.85", "3.88", "0.3434", 6),
    @("202801", "南方全球精选配置(QDII-FOF)", "18.00", "28.82", "1.64", "0.2952", 4),
    @("011315", "永赢港股通优质成长一年持有期混合型证券投资基金", "4.07", "90.36", "6.18", "0.2515", 3),
    @("457001", "国富亚洲机会股票 (QDII)", "5.93", "77.36", "3.00", "0.1779", 4),
    @("011905", "安信价值启航混合型证券投资基金A", "4.48", "93.14", "3.87", "0.1734", 6),
    @("008477", "安信价值驱动三年持有期混合", "2.64", "91.75", "3.76", "0.0993", 6),
    @("011117", "富国沪港深业绩驱动混合C", "2.39", "74.38", "2.95", "0.0705", 6),
    @("161229", "国投瑞银中国价值发现股票QDII-LOF", "1.47", "92.83", "4.60", "0.0676", 7),
    @("519139", "海富通沪港深灵活配置混合", "1.32", "94.37", "4.21", "0.0556", 7),
    @("011969", "建信港股通精选混合A", "1.01", "57.64", "5.33", "0.0538", 4),
    @("920011", "中金安心回报灵活配置混合A", "4.29", "54.53", "1.06", "0.0455", 8),
    @("008253", "华宝致远混合（QDII）A", "0.70", "85.00", "3.74", "0.0262", 8),
    @("010105", "景顺长城消费精选混合C", "0.51", "90.93", "4.42", "0.0225", 7),
    @("920921", "中金安心回报灵活配置混合C", "1.68", "54.53", "1.06", "0.0178", 8),
    @("011970", "建信港股通精选混合C", "0.33", "57.64", "5.33", "0.0176", 4),
    @("011906", "安信价值启航混合型证券投资基金C", "0.42", "93.14", "3.87", "0.0163", 6),
    @("000927", "博时大中华亚太精选股票(QDII) - 美元现汇", "0.32", "92.94", "4.31", "0.0138", 7),
    @("050015", "博时大中华亚太精选股票(QDII) -人民币", "0.32", "92.94", "4.31", "0.0138", 7),
    @("007505", "华夏中证AH经济蓝筹股票指数A", "1.09", "93.35", "1.14", "0.0124", 10),
    @("004098", "前海开源港股通股息率50强股票", "0.34", "88.92", "3.12", "0.0106", 6),
    @("007506", "华夏中证AH经济蓝筹股票指数C", "0.65", "93.35", "1.14", "0.0074", 10),
    @("008254", "华宝致远混合（QDII）C", "0.12", "85.00", "3.74", "0.0045", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = ($r - 2)

    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[5]

    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Row 35 is a special case: the fund had no reported "基金规模" (column D
# is blank) and its 持有市值 rounds to a genuine numeric 0 (not text "0").
$newSheet.Cells.Item(35, 1).Value = 33
$newSheet.Range("B35").NumberFormat = "@"
$newSheet.Cells.Item(35, 2).Value = "007234"
$newSheet.Cells.Item(35, 3).Value = "博时优势企业3年封闭运作灵活配置混合（LOF）C"
$newSheet.Range("D35").ClearContents()
$newSheet.Range("E35").NumberFormat = "@"
$newSheet.Cells.Item(35, 5).Value = "82.86"
$newSheet.Range("F35").NumberFormat = "@"
$newSheet.Cells.Item(35, 6).Value = "4.38"
$newSheet.Cells.Item(35, 7).Value = 0
$newSheet.Cells.Item(35, 8).Value = 6

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert a new row for 2022-Q1 above
#    the existing 2021-Q4 row, and renumber the index column (A).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Row-insert leaves the freshly inserted row unformatted, so copy the
# index-column style back in from the row just below (which carries the
# original row-2 formatting, shifted down to row 3 by the insert).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 34
$totalSheet.Cells.Item(2, 4).Value = 10.98

for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

$totalSheet.Range("A1").Select()
